$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Get-ParagraphContaining {
    param($doc, [string]$needle)
    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Text.Contains($needle)) {
            return $para
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Remove the stray <w:proofErr w:type="gramStart"/> / gramEnd marker that
#    wraps "canvas.width" (the spellStart/spellEnd pair stays untouched).
# ---------------------------------------------------------------------------
$p1 = (Get-ParagraphContaining $d "canvas.width").Range
$p1xml = '<w:p w14:paraId="5D8A680C" w14:textId="5D7795C3" w:rsidR="00A61938" w:rsidRPr="00A82067" w:rsidRDefault="00A61938" w:rsidP="00A61938" ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="00A82067"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Restrict that the text we add on the meme cannot be longer than the </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r w:rsidRPr="00A82067"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>canvas.width</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r w:rsidR="008430E4" w:rsidRPr="00A82067"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> - Done</w:t></w:r>' +
  '</w:p>'
[void]$p1.InsertXML($p1xml)

# ---------------------------------------------------------------------------
# 2) Remove the stray <w:proofErr w:type="gramStart"/> / gramEnd marker that
#    wraps "it's" (the spellStart/spellEnd pair stays untouched).
# ---------------------------------------------------------------------------
$p2 = (Get-ParagraphContaining $d "width and height, good proportions").Range
$p2xml = '<w:p w14:paraId="552C85E5" w14:textId="0F18AAB3" w:rsidR="00A61938" w:rsidRPr="0055222B" w:rsidRDefault="00635F61" w:rsidP="00A61938" ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="0055222B"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Image in the canvas should be in </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r w:rsidRPr="0055222B"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>it’s</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r w:rsidRPr="0055222B"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> width and height, good proportions</w:t></w:r>' +
  '<w:r w:rsidR="007A6355" w:rsidRPr="0055222B"><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> - Done</w:t></w:r>' +
  '</w:p>'
[void]$p2.InsertXML($p2xml)

# ---------------------------------------------------------------------------
# 3) Split the last bullet ("When clicking on search keyword ...") so the
#    trailing " - " becomes its own run using an en-dash plus a space, and
#    append a brand new bullet "Continue the bonuses - ".
# ---------------------------------------------------------------------------
$last = (Get-ParagraphContaining $d "When clicking on search keyword").Range
$lastXml = '<w:p ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">When clicking on search keyword make it bigger and can see it without refreshing the page </w:t></w:r>' +
  '<w:r><w:t>–</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '</w:p>' +
  '<w:p>' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Continue the bonuses - </w:t></w:r>' +
  '</w:p>'
[void]$last.InsertXML($lastXml)
